$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block (rows 77..173) records weekly "Ajo" (garlic) price quotes.
# A new week's record was inserted at the top of that block (row 77), which
# pushes every existing record down by one row; the record that used to be
# last (old row 173) becomes the new last row (174).
#
# Column D (Fecha) and columns I..P (Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Unidad de comercializacion,
# Origen, Precio $/Kg) are the per-record fields that move down by one row.
# Columns A,B,C,E,F,G,H,Q,R are constant across the whole block (Mercado,
# Region, Codreg, Categoria, Variedad, Kg o Unidades, Clasificacion), so they
# don't need to move - row 174 just gets the same constants.

$lastOldRow = 173
$firstOldRow = 76
$newLastRow = 174
$firstNewRow = 77

# Shift column D (Fecha) down by one row.
$srcD = $ws.Range("D$firstOldRow`:D$lastOldRow").Value2
$ws.Range("D$firstNewRow`:D$newLastRow").Value2 = $srcD

# Shift the contiguous block I:P (Calidad..Precio $/Kg) down by one row.
$srcIP = $ws.Range("I$firstOldRow`:P$lastOldRow").Value2
$ws.Range("I$firstNewRow`:P$newLastRow").Value2 = $srcIP

# New row 174 keeps the same constant values as the rest of the block for
# the columns that never change.
$ws.Cells.Item($newLastRow,1).Value2  = $ws.Cells.Item($lastOldRow,1).Value2   # A Mercado ID
$ws.Cells.Item($newLastRow,2).Value2  = $ws.Cells.Item($lastOldRow,2).Value2   # B Mercado
$ws.Cells.Item($newLastRow,3).Value2  = $ws.Cells.Item($lastOldRow,3).Value2   # C Region
$ws.Cells.Item($newLastRow,5).Value2  = $ws.Cells.Item($lastOldRow,5).Value2   # E Codreg
$ws.Cells.Item($newLastRow,6).Value2  = $ws.Cells.Item($lastOldRow,6).Value2   # F Categoria ID
$ws.Cells.Item($newLastRow,7).Value2  = $ws.Cells.Item($lastOldRow,7).Value2   # G Categoria
$ws.Cells.Item($newLastRow,8).Value2  = $ws.Cells.Item($lastOldRow,8).Value2   # H Variedad
$ws.Cells.Item($newLastRow,17).Value2 = $ws.Cells.Item($lastOldRow,17).Value2  # Q Kg o Unidades
$ws.Cells.Item($newLastRow,18).Value2 = $ws.Cells.Item($lastOldRow,18).Value2  # R Clasificacion

# Match the date-formatted style of column D for the new row.
$ws.Cells.Item($newLastRow,4).NumberFormat = $ws.Cells.Item($lastOldRow,4).NumberFormat
